# ajout/deplacement de la colonne commentaire
#
# Inserts a new "Commentaire" column before the current "Profil" column
# (i.e. before column F), pushing the existing F:M columns to G:N, fills
# in the header/value for the new column, restores the column width,
# extends the AutoFilter (and the _FilterDatabase defined name) to cover
# the new column, and restores the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank column at F, shifting Profil..Ordre (F..M) to G..N
$ws.Columns("F:F").Insert()

# 2) Populate the new column's header (row1) and sample value (row2)
$ws.Range("F1").Value = "Commentaire"
$ws.Range("F2").Value = "test"

# 3) Restore the column width for the new column (closest attainable value)
$ws.Columns("F:F").ColumnWidth = 32.1

# 4) Extend the AutoFilter range from C1:J2 to C1:K2 (now includes col K,
#    previously col J, since everything shifted right by one column)
$ws.AutoFilterMode = $false
$ws.Range("C1:K2").AutoFilter() | Out-Null

# 5) Update the _FilterDatabase defined name to match the new AutoFilter range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Exigences!_FilterDatabase") {
        $n.RefersTo = "=Exigences!`$C`$1:`$K`$1"
    }
}

# 6) Restore the active cell selection on the sheet
$ws.Range("F6").Select() | Out-Null
